$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap the B1/C1 labels ---
$ws.Range("B1").Value = "MLBSO00"
$ws.Range("C1").Value = "LNBSF00"

# --- Data rows 2-17: swap B/C values for each existing row ---
$data = @(
    @(45734, 806.651, 790.4),
    @(45733, 810.465, 795.08),
    @(45730, 810.465, 792.8440000000001),
    @(45729, 810.465, 803.816),
    @(45728, 810.465, 802.724),
    @(45735, 806.651, 808.9640000000001),
    @(45736, 806.651, 823.3680000000001),
    @(45737, 806.651, 823.9400000000001),
    @(45740, 806.651, 788.6319999999999),
    @(45741, 806.651, 785.928),
    @(45742, 806.651, 773.9160000000001),
    @(45734, 806.651, 790.4),
    @(45743, 806.651, 771.1079999999999),
    @(45744, 800.9299999999999, 760.188),
    @(45748, 800.9299999999999, 749.736),
    @(45749, 800.9299999999999, 764.244),
    @(45750, 797.116, 753.74)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# --- Row 18 is new: fill in column A (date) too ---
$ws.Cells.Item(18, 1).Value = 45750

# --- Row 17's date style changes from the "date only" format to the
#     "date time" format used by the other data rows, and the new row 18
#     now carries the "date only" format that used to sit on row 17 ---
$ws.Range("A17").NumberFormat = $ws.Range("A16").NumberFormat
$ws.Range("A18").NumberFormat = "YYYY-MM-DD"
